$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025-03-20")

$values = @{
    2  = @(65, 220)
    3  = @(84, 178)
    4  = @(96, 277)
    5  = @(100, 218)
    6  = @(249, 625)
    7  = @(995, 995)
    8  = @(97, 257)
    9  = @(100, 253)
    10 = @(105, 246)
    11 = @(107, 327)
    12 = @(233, 858)
    13 = @($null, 346)
    14 = @(148, 274)
    15 = @(149, 238)
    16 = @(155, 269)
    17 = @(261, 261)
    18 = @(299, 716)
    19 = @(608, 747)
    20 = @(615, 879)
    21 = @(662, 869)
    22 = @($null, 1214)
    23 = @(967, 1641)
    24 = @(1658, 2488)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $cVal = $pair[0]
    $dVal = $pair[1]
    if ($null -ne $cVal) {
        $ws.Cells.Item($row, 3).Value = $cVal
    }
    $ws.Cells.Item($row, 4).Value = $dVal
}
